# Auto-applies the gh-pages data-refresh commit: updated "want-to-go" counts (column F)
# across sheets, plus a new exhibition row ("杭州·白日梦次元动漫嘉年华") inserted into the
# "展览" (Exhibitions) sheet and the "全部类型" (All types) combined sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item(1)

# Updated "want to go" counts for existing rows
$ws1.Range("F2").Value = 23
$ws1.Range("F3").Value = 8390
$ws1.Range("F6").Value = 308
$ws1.Range("F8").Value = 644
$ws1.Range("F9").Value = 118
$ws1.Range("F12").Value = 886
$ws1.Range("F13").Value = 3622
$ws1.Range("F14").Value = 252
$ws1.Range("F15").Value = 140
$ws1.Range("F16").Value = 782
$ws1.Range("F17").Value = 768
$ws1.Range("F18").Value = 52
$ws1.Range("F19").Value = 488
$ws1.Range("F22").Value = 1026
$ws1.Range("F23").Value = 1350
$ws1.Range("F24").Value = 415
$ws1.Range("F25").Value = 289
$ws1.Range("F26").Value = 142
$ws1.Range("F27").Value = 147
$ws1.Range("F28").Value = 328
$ws1.Range("F29").Value = 50
$ws1.Range("F30").Value = 1017
$ws1.Range("F32").Value = 511
$ws1.Range("F33").Value = 635
$ws1.Range("F34").Value = 41
$ws1.Range("F36").Value = 67
$ws1.Range("F37").Value = 27
$ws1.Range("F39").Value = 145

# New row 40: 杭州·白日梦次元动漫嘉年华 (2024-04-20)
$ws1.Range("A39:I39").Copy()
$ws1.Range("A40:I40").PasteSpecial(-4122)
$ws1.Cells.Item(40, 1).Value = 39
$ws1.Cells.Item(40, 2).Value = "'2024-04-20"
$ws1.Cells.Item(40, 3).Value = "杭州·白日梦次元动漫嘉年华"
$ws1.Cells.Item(40, 4).Value = "黄姑山路51-4号 0101park"
$ws1.Cells.Item(40, 5).Value = "2024.04.20 10:00-04.21 18:00"
$ws1.Cells.Item(40, 6).Value = 1
$ws1.Cells.Item(40, 7).Value = 68
$ws1.Cells.Item(40, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81634"
$ws1.Cells.Item(40, 9).Value = "//i0.hdslb.com/bfs/openplatform/202402/n65ZLevi1706777788165.jpeg"

# --- Sheet 2: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 23

# --- Sheet 4: 全部类型 (All types, combined) ---
$ws4 = $wb.Worksheets.Item(4)

# Updated "want to go" counts for existing rows
$ws4.Range("F3").Value = 23
$ws4.Range("F4").Value = 8390
$ws4.Range("F7").Value = 308
$ws4.Range("F9").Value = 644
$ws4.Range("F10").Value = 118
$ws4.Range("F13").Value = 886
$ws4.Range("F14").Value = 23
$ws4.Range("F15").Value = 3622
$ws4.Range("F16").Value = 252
$ws4.Range("F17").Value = 140
$ws4.Range("F19").Value = 782
$ws4.Range("F20").Value = 768
$ws4.Range("F22").Value = 52
$ws4.Range("F23").Value = 488
$ws4.Range("F27").Value = 1026
$ws4.Range("F28").Value = 1350
$ws4.Range("F29").Value = 415
$ws4.Range("F30").Value = 289
$ws4.Range("F31").Value = 142
$ws4.Range("F32").Value = 147
$ws4.Range("F34").Value = 328
$ws4.Range("F35").Value = 50
$ws4.Range("F36").Value = 1017
$ws4.Range("F38").Value = 511
$ws4.Range("F39").Value = 635
$ws4.Range("F40").Value = 41
$ws4.Range("F42").Value = 67
$ws4.Range("F43").Value = 27
$ws4.Range("F45").Value = 145

# Insert new row 46: 杭州·白日梦次元动漫嘉年华 (2024-04-20), pushing the two rows after it down
$ws4.Rows.Item(46).Insert()

# Re-index the "#" column for the two rows that shifted down
$ws4.Cells.Item(47, 1).Value = 46
$ws4.Cells.Item(48, 1).Value = 47

# Give the freshly inserted row the same look (border/bold/centered) as the other data rows
$ws4.Range("A45:I45").Copy()
$ws4.Range("A46:I46").PasteSpecial(-4122)

$ws4.Cells.Item(46, 1).Value = 45
$ws4.Cells.Item(46, 2).Value = "'2024-04-20"
$ws4.Cells.Item(46, 3).Value = "杭州·白日梦次元动漫嘉年华"
$ws4.Cells.Item(46, 4).Value = "黄姑山路51-4号 0101park"
$ws4.Cells.Item(46, 5).Value = "2024.04.20 10:00-04.21 18:00"
$ws4.Cells.Item(46, 6).Value = 1
$ws4.Cells.Item(46, 7).Value = 68
$ws4.Cells.Item(46, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81634"
$ws4.Cells.Item(46, 9).Value = "//i0.hdslb.com/bfs/openplatform/202402/n65ZLevi1706777788165.jpeg"

